$wb = $excel.ActiveWorkbook
$ws1 = $wb.Sheets.Item("Moorings")
$ws2 = $wb.Sheets.Item("Asset_Cal_Info")

# Sheet1 (Moorings): set Mooring OOIBARCODE value for A2
$ws1.Range("A2").ClearFormats()
$ws1.Range("A2").Value = "OL000288"

# Sheet2 (Asset_Cal_Info): fill in Mooring OOIBARCODE (B) and Sensor OOIBARCODE (E) columns
$ws2.Range("B3").ClearFormats()
$ws2.Range("B3").Value = "OL000288"
$ws2.Range("E3").ClearFormats()
$ws2.Range("E3").Value = "A00285"
$ws2.Range("B5").ClearFormats()
$ws2.Range("B5").Value = "OL000288"
$ws2.Range("E5").ClearFormats()
$ws2.Range("E5").Value = "OL000289"
$ws2.Range("B6").ClearFormats()
$ws2.Range("B6").Value = "OL000288"
$ws2.Range("E6").ClearFormats()
$ws2.Range("E6").Value = "OL000289"
$ws2.Range("B8").ClearFormats()
$ws2.Range("B8").Value = "OL000288"
$ws2.Range("E8").ClearFormats()
$ws2.Range("E8").Value = "N00051"
$ws2.Range("B9").ClearFormats()
$ws2.Range("B9").Value = "OL000288"
$ws2.Range("E9").ClearFormats()
$ws2.Range("E9").Value = "N00051"
$ws2.Range("B10").ClearFormats()
$ws2.Range("B10").Value = "OL000288"
$ws2.Range("E10").ClearFormats()
$ws2.Range("E10").Value = "N00051"
$ws2.Range("B11").ClearFormats()
$ws2.Range("B11").Value = "OL000288"
$ws2.Range("E11").ClearFormats()
$ws2.Range("E11").Value = "N00051"
$ws2.Range("B12").ClearFormats()
$ws2.Range("B12").Value = "OL000288"
$ws2.Range("E12").ClearFormats()
$ws2.Range("E12").Value = "N00051"
$ws2.Range("B13").ClearFormats()
$ws2.Range("B13").Value = "OL000288"
$ws2.Range("E13").ClearFormats()
$ws2.Range("E13").Value = "N00051"
$ws2.Range("B14").ClearFormats()
$ws2.Range("B14").Value = "OL000288"
$ws2.Range("E14").ClearFormats()
$ws2.Range("E14").Value = "N00051"
$ws2.Range("B15").ClearFormats()
$ws2.Range("B15").Value = "OL000288"
$ws2.Range("E15").ClearFormats()
$ws2.Range("E15").Value = "N00051"
$ws2.Range("B17").ClearFormats()
$ws2.Range("B17").Value = "OL000288"
$ws2.Range("E17").ClearFormats()
$ws2.Range("E17").Value = "N00052"
$ws2.Range("B18").ClearFormats()
$ws2.Range("B18").Value = "OL000288"
$ws2.Range("E18").ClearFormats()
$ws2.Range("E18").Value = "N00052"
$ws2.Range("B19").ClearFormats()
$ws2.Range("B19").Value = "OL000288"
$ws2.Range("E19").ClearFormats()
$ws2.Range("E19").Value = "N00052"
$ws2.Range("B20").ClearFormats()
$ws2.Range("B20").Value = "OL000288"
$ws2.Range("E20").ClearFormats()
$ws2.Range("E20").Value = "N00052"
$ws2.Range("B21").ClearFormats()
$ws2.Range("B21").Value = "OL000288"
$ws2.Range("E21").ClearFormats()
$ws2.Range("E21").Value = "N00052"
$ws2.Range("B22").ClearFormats()
$ws2.Range("B22").Value = "OL000288"
$ws2.Range("E22").ClearFormats()
$ws2.Range("E22").Value = "N00052"
$ws2.Range("B23").ClearFormats()
$ws2.Range("B23").Value = "OL000288"
$ws2.Range("E23").ClearFormats()
$ws2.Range("E23").Value = "N00052"
$ws2.Range("B24").ClearFormats()
$ws2.Range("B24").Value = "OL000288"
$ws2.Range("E24").ClearFormats()
$ws2.Range("E24").Value = "N00052"
$ws2.Range("B25").ClearFormats()
$ws2.Range("B25").Value = "OL000288"
$ws2.Range("E25").ClearFormats()
$ws2.Range("E25").Value = "N00052"
$ws2.Range("B26").ClearFormats()
$ws2.Range("B26").Value = "OL000288"
$ws2.Range("E26").ClearFormats()
$ws2.Range("E26").Value = "N00052"
$ws2.Range("B28").ClearFormats()
$ws2.Range("B28").Value = "OL000288"
$ws2.Range("E28").ClearFormats()
$ws2.Range("E28").Value = "N00054"
$ws2.Range("B29").ClearFormats()
$ws2.Range("B29").Value = "OL000288"
$ws2.Range("E29").ClearFormats()
$ws2.Range("E29").Value = "N00054"
$ws2.Range("B31").ClearFormats()
$ws2.Range("B31").Value = "OL000288"
$ws2.Range("E31").ClearFormats()
$ws2.Range("E31").Value = "N00055"
$ws2.Range("B32").ClearFormats()
$ws2.Range("B32").Value = "OL000288"
$ws2.Range("E32").ClearFormats()
$ws2.Range("E32").Value = "N00055"
$ws2.Range("B34").ClearFormats()
$ws2.Range("B34").Value = "OL000288"
$ws2.Range("E34").ClearFormats()
$ws2.Range("E34").Value = "OL000242"

# Restore selections (without changing the active sheet/tab)
$ws1.Range("A2").Select()
$ws2.Range("F32").Select()
